$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185574173927307
$ws.Range("B1").Value = 2.133893489837646
$ws.Range("C1").Value = 3.780712366104126
$ws.Range("D1").Value = 3.268907070159912
$ws.Range("E1").Value = 1.143660306930542
